# =====================================================================
# Leave Card update — "Update Leave Card 4/12/2023 4:43 PM"
# Fills in the monthly leave-accrual rows for Nov 2022 through Jan 2030,
# records the Dec 2022 Forced-Leave / Mar 2023 Special-Leave remarks,
# inserts a new blank table row before the final row, and nudges the
# saved view state to where the author last had it.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
# --- Fill in leave-card rows 291-378 (monthly entries extending the table) ---
$ws.Cells.Item(291,1).Value = 44866
$ws.Cells.Item(291,3).Value = 1.25
$ws.Cells.Item(292,1).Value = 44896
$ws.Cells.Item(292,2).Value = "FL(5-0-0)"
$ws.Cells.Item(292,3).Value = 1.25
$ws.Cells.Item(292,4).Value = 5
$ws.Cells.Item(293,1).Value = "'2023"
$ws.Range("A280").Copy()
$ws.Range("A293").PasteSpecial(-4122)
$ws.Cells.Item(294,1).Value = 44927
$ws.Cells.Item(294,3).Value = 1.25
$ws.Cells.Item(295,1).Value = 44958
$ws.Cells.Item(295,3).Value = 1.25
$ws.Cells.Item(296,1).Value = 44986
$ws.Cells.Item(296,2).Value = "SP(1-0-0)"
$ws.Cells.Item(296,3).Value = 1.25
$ws.Range("K23").Copy()
$ws.Range("K296").PasteSpecial(-4122)
$ws.Cells.Item(296,11).Value = 45005
$ws.Cells.Item(297,1).Value = 45017
$ws.Cells.Item(298,1).Value = 45047
$ws.Cells.Item(299,1).Value = 45078
$ws.Cells.Item(300,1).Value = 45108
$ws.Cells.Item(301,1).Value = 45139
$ws.Cells.Item(302,1).Value = 45170
$ws.Cells.Item(303,1).Value = 45200
$ws.Cells.Item(304,1).Value = 45231
$ws.Cells.Item(305,1).Value = 45261
$ws.Cells.Item(306,1).Value = 45292
$ws.Cells.Item(307,1).Value = 45323
$ws.Cells.Item(308,1).Value = 45352
$ws.Cells.Item(309,1).Value = 45383
$ws.Cells.Item(310,1).Value = 45413
$ws.Cells.Item(311,1).Value = 45444
$ws.Cells.Item(312,1).Value = 45474
$ws.Cells.Item(313,1).Value = 45505
$ws.Cells.Item(314,1).Value = 45536
$ws.Cells.Item(315,1).Value = 45566
$ws.Cells.Item(316,1).Value = 45597
$ws.Cells.Item(317,1).Value = 45627
$ws.Cells.Item(318,1).Value = 45658
$ws.Cells.Item(319,1).Value = 45689
$ws.Cells.Item(320,1).Value = 45717
$ws.Cells.Item(321,1).Value = 45748
$ws.Cells.Item(322,1).Value = 45778
$ws.Cells.Item(323,1).Value = 45809
$ws.Cells.Item(324,1).Value = 45839
$ws.Cells.Item(325,1).Value = 45870
$ws.Cells.Item(326,1).Value = 45901
$ws.Cells.Item(327,1).Value = 45931
$ws.Cells.Item(328,1).Value = 45962
$ws.Cells.Item(329,1).Value = 45992
$ws.Cells.Item(330,1).Value = 46023
$ws.Cells.Item(331,1).Value = 46054
$ws.Cells.Item(332,1).Value = 46082
$ws.Cells.Item(333,1).Value = 46113
$ws.Cells.Item(334,1).Value = 46143
$ws.Cells.Item(335,1).Value = 46174
$ws.Cells.Item(336,1).Value = 46204
$ws.Cells.Item(337,1).Value = 46235
$ws.Cells.Item(338,1).Value = 46266
$ws.Cells.Item(339,1).Value = 46296
$ws.Cells.Item(340,1).Value = 46327
$ws.Cells.Item(341,1).Value = 46357
$ws.Cells.Item(342,1).Value = 46388
$ws.Cells.Item(343,1).Value = 46419
$ws.Cells.Item(344,1).Value = 46447
$ws.Cells.Item(345,1).Value = 46478
$ws.Cells.Item(346,1).Value = 46508
$ws.Cells.Item(347,1).Value = 46539
$ws.Cells.Item(348,1).Value = 46569
$ws.Cells.Item(349,1).Value = 46600
$ws.Cells.Item(350,1).Value = 46631
$ws.Cells.Item(351,1).Value = 46661
$ws.Cells.Item(352,1).Value = 46692
$ws.Cells.Item(353,1).Value = 46722
$ws.Cells.Item(354,1).Value = 46753
$ws.Cells.Item(355,1).Value = 46784
$ws.Cells.Item(356,1).Value = 46813
$ws.Cells.Item(357,1).Value = 46844
$ws.Cells.Item(358,1).Value = 46874
$ws.Cells.Item(359,1).Value = 46905
$ws.Cells.Item(360,1).Value = 46935
$ws.Cells.Item(361,1).Value = 46966
$ws.Cells.Item(362,1).Value = 46997
$ws.Cells.Item(363,1).Value = 47027
$ws.Cells.Item(364,1).Value = 47058
$ws.Cells.Item(365,1).Value = 47088
$ws.Cells.Item(366,1).Value = 47119
$ws.Cells.Item(367,1).Value = 47150
$ws.Cells.Item(368,1).Value = 47178
$ws.Cells.Item(369,1).Value = 47209
$ws.Cells.Item(370,1).Value = 47239
$ws.Cells.Item(371,1).Value = 47270
$ws.Cells.Item(372,1).Value = 47300
$ws.Cells.Item(373,1).Value = 47331
$ws.Cells.Item(374,1).Value = 47362
$ws.Cells.Item(375,1).Value = 47392
$ws.Cells.Item(376,1).Value = 47423
$ws.Cells.Item(377,1).Value = 47453
$ws.Cells.Item(378,1).Value = 47484
# --- Insert a new (blank) row into Table1 right before the final row ---
# Table1 currently spans A8:K749; row 749 carries the table's distinctive
# "final row" bottom-border formatting. We want that formatting to end up
# one row lower (on the new last row, 750) while row 749 becomes an
# ordinary data row (matching row 748's look), and the table + sheet
# dimension both grow to K750.
$lo = $ws.ListObjects.Item("Table1")

# 1) Snapshot the current "final row" formatting before we touch anything.
$ws.Range("A749:K749").Copy()

# 2) Grow the table by one row.
$lo.Resize($ws.Range("A8:K750"))

# 3) Apply the snapshot to the new final row.
$ws.Range("A750:K750").PasteSpecial(-4122)

# 4) Row 749 should now look like a normal data row, like row 748.
$ws.Range("A748:K748").Copy()
$ws.Range("A749:K749").PasteSpecial(-4122)

# 5) Restore the calculated-column formulas (PasteSpecial formats only,
#    but Resize/Copy can leave the "EARNED " helper formula stale).
$ws.Range("G749").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G750").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Restore the saved view state (split position + last selections) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 225
$ws.Range("O27").Select()
$ws.Range("E234").Select()
